$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2.5
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("W2").Value = 9
$ws.Range("AC2").Value = 17

$ws.Range("G3").Value = 1.91
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.33
$ws.Range("K3").Value = 2
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("X3").Value = 8
$ws.Range("AA3").Value = 19
$ws.Range("AC3").Value = 7
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 501
$ws.Range("AH3").Value = 9.5
$ws.Range("AI3").Value = 21
$ws.Range("AK3").Value = 51
$ws.Range("AM3").Value = 51
$ws.Range("AP3").Value = 26
$ws.Range("AU3").Value = 9
$ws.Range("AX3").Value = 26
$ws.Range("BB3").Value = 351

$ws.Range("G4").Value = 2.5
$ws.Range("I4").Value = 2.6
$ws.Range("J4").Value = 3
$ws.Range("L4").Value = 3.1
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 2.3
$ws.Range("W4").Value = 12
$ws.Range("AA4").Value = 17
$ws.Range("AL4").Value = 19
$ws.Range("AM4").Value = 23
$ws.Range("AN4").Value = 4.75
$ws.Range("AW4").Value = 5

$ws.Range("G5").Value = 1.95
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.63
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 4.5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("U5").Value = 1.95
$ws.Range("V5").Value = 1.8
$ws.Range("W5").Value = 6.5
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 17
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 8.5
$ws.Range("AD5").Value = 6
$ws.Range("AG5").Value = 301
$ws.Range("AH5").Value = 11
$ws.Range("AI5").Value = 21
$ws.Range("AK5").Value = 41
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 11
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 41
$ws.Range("AT5").Value = 2.63
$ws.Range("AU5").Value = 8.5
$ws.Range("AW5").Value = 5.5
$ws.Range("AY5").Value = 34
$ws.Range("BB5").Value = 251

$ws.Range("G6").Value = 3.9
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 1.75
$ws.Range("J6").Value = 4
$ws.Range("L6").Value = 2.25
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 23
$ws.Range("O6").Value = 1.11
$ws.Range("P6").Value = 6.5
$ws.Range("S6").Value = 1.22
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = 1.44
$ws.Range("V6").Value = 2.63
$ws.Range("X6").Value = 26
$ws.Range("Y6").Value = 13
$ws.Range("AA6").Value = 26
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 23
$ws.Range("AE6").Value = 12
$ws.Range("AF6").Value = 29
$ws.Range("AG6").Value = 81
$ws.Range("AH6").Value = 13
$ws.Range("AI6").Value = 12
$ws.Range("AJ6").Value = 9
$ws.Range("AK6").Value = 17
$ws.Range("AL6").Value = 12
$ws.Range("AN6").Value = 6.5
$ws.Range("AO6").Value = 19
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 51
$ws.Range("AS6").Value = 101
$ws.Range("AT6").Value = 4
$ws.Range("AU6").Value = 7
$ws.Range("AW6").Value = 4.33
$ws.Range("AX6").Value = 9
$ws.Range("AZ6").Value = 26
$ws.Range("BA6").Value = 34
$ws.Range("BC6").Value = 251

$ws.Range("O7").Value = 1.17
$ws.Range("P7").Value = 5
$ws.Range("Q7").Value = 1.57
$ws.Range("R7").Value = 2.35
$ws.Range("U7").Value = 1.5
$ws.Range("V7").Value = 2.5
$ws.Range("AG7").Value = 101
$ws.Range("AI7").Value = 19
$ws.Range("AR7").Value = 41

$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("Q10").Value = 2.08
$ws.Range("R10").Value = 1.73
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("AC10").Value = 9
$ws.Range("AT10").Value = 2.63
